$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (new volume number / week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/13/2025  Through  1/19/2025"

# --- Cells whose type switches between blank-placeholder text and a real number ---
# (copy number format from a same-style neighbor, then set the value, so the
#  underlying style index is reused instead of minting a new one)
$ws.Range("D29").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 3

$ws.Range("D29").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = 2

$ws.Range("D29").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D31").Value = 1

$ws.Range("E29").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E31").Value = -100

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D33").PasteSpecial(-4122)

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E33").PasteSpecial(-4122)

# --- Plain numeric refresh for the rest of the weekly crime-stat table ---
# Row 14
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 5
$ws.Range("H14").Value = -80
$ws.Range("J14").Value = 4

# Row 15
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 25
$ws.Range("F15").Value = 18
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 14
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = 75
$ws.Range("L15").Value = -30
$ws.Range("M15").Value = 27.272727272727
$ws.Range("N15").Value = -51.724137931034

# Row 16
$ws.Range("C16").Value = 26
$ws.Range("D16").Value = 42
$ws.Range("E16").Value = -38.095238095238
$ws.Range("F16").Value = 107
$ws.Range("G16").Value = 192
$ws.Range("H16").Value = -44.270833333333
$ws.Range("I16").Value = 71
$ws.Range("J16").Value = 121
$ws.Range("K16").Value = -41.322314049586
$ws.Range("L16").Value = -41.322314049586
$ws.Range("M16").Value = -62.032085561497
$ws.Range("N16").Value = -92.146017699115

# Row 17
$ws.Range("C17").Value = 65
$ws.Range("D17").Value = 72
$ws.Range("E17").Value = -9.722222222222
$ws.Range("F17").Value = 252
$ws.Range("G17").Value = 311
$ws.Range("H17").Value = -18.971061093247
$ws.Range("I17").Value = 169
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = -15.5
$ws.Range("L17").Value = -15.075376884422
$ws.Range("M17").Value = 6.962025316455
$ws.Range("N17").Value = -58.780487804878

# Row 18
$ws.Range("C18").Value = 23
$ws.Range("D18").Value = 38
$ws.Range("E18").Value = -39.473684210526
$ws.Range("F18").Value = 108
$ws.Range("G18").Value = 155
$ws.Range("H18").Value = -30.322580645161
$ws.Range("I18").Value = 68
$ws.Range("J18").Value = 104
$ws.Range("K18").Value = -34.615384615384
$ws.Range("L18").Value = -55.555555555555
$ws.Range("M18").Value = -53.424657534246
$ws.Range("N18").Value = -88.235294117647

# Row 19
$ws.Range("C19").Value = 71
$ws.Range("D19").Value = 80
$ws.Range("E19").Value = -11.25
$ws.Range("F19").Value = 263
$ws.Range("G19").Value = 358
$ws.Range("H19").Value = -26.536312849162
$ws.Range("I19").Value = 180
$ws.Range("J19").Value = 240
$ws.Range("K19").Value = -25
$ws.Range("L19").Value = -37.062937062937
$ws.Range("M19").Value = 7.142857142857
$ws.Range("N19").Value = -41.558441558441

# Row 20
$ws.Range("C20").Value = 16
$ws.Range("D20").Value = 33
$ws.Range("E20").Value = -51.515151515151
$ws.Range("F20").Value = 86
$ws.Range("G20").Value = 130
$ws.Range("H20").Value = -33.846153846153
$ws.Range("I20").Value = 58
$ws.Range("J20").Value = 96
$ws.Range("K20").Value = -39.583333333333
$ws.Range("L20").Value = -30.12048192771
$ws.Range("M20").Value = -12.121212121212
$ws.Range("N20").Value = -89.199255121042

# Row 21
$ws.Range("C21").Value = 206
$ws.Range("D21").Value = 271
$ws.Range("E21").Value = -23.985239852398
$ws.Range("F21").Value = 835
$ws.Range("G21").Value = 1163
$ws.Range("H21").Value = -28.202923473774
$ws.Range("I21").Value = 560
$ws.Range("J21").Value = 773
$ws.Range("K21").Value = -27.554980595084
$ws.Range("L21").Value = -35.409457900807
$ws.Range("M21").Value = -24.42645074224
$ws.Range("N21").Value = -79.964221824686

# Row 22
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 12
$ws.Range("E22").Value = -58.333333333333
$ws.Range("F22").Value = 19
$ws.Range("G22").Value = 32
$ws.Range("H22").Value = -40.625
$ws.Range("I22").Value = 13
$ws.Range("J22").Value = 21
$ws.Range("K22").Value = -38.095238095238
$ws.Range("L22").Value = -40.90909090909
$ws.Range("M22").Value = -23.529411764705

# Row 23
$ws.Range("D23").Value = 30
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 88
$ws.Range("G23").Value = 117
$ws.Range("H23").Value = -24.786324786324
$ws.Range("I23").Value = 60
$ws.Range("J23").Value = 78
$ws.Range("K23").Value = -23.076923076923
$ws.Range("L23").Value = -28.571428571428
$ws.Range("M23").Value = 30.434782608695

# Row 24
$ws.Range("C24").Value = 214
$ws.Range("D24").Value = 165
$ws.Range("E24").Value = 29.696969696969
$ws.Range("F24").Value = 858
$ws.Range("G24").Value = 797
$ws.Range("H24").Value = 7.653701380175
$ws.Range("I24").Value = 556
$ws.Range("J24").Value = 521
$ws.Range("K24").Value = 6.717850287907
$ws.Range("L24").Value = -9.001636661211
$ws.Range("M24").Value = 10.09900990099

# Row 25
$ws.Range("C25").Value = 65
$ws.Range("D25").Value = 70
$ws.Range("E25").Value = -7.142857142857
$ws.Range("F25").Value = 302
$ws.Range("G25").Value = 290
$ws.Range("H25").Value = 4.137931034482
$ws.Range("I25").Value = 193
$ws.Range("J25").Value = 190
$ws.Range("K25").Value = 1.578947368421
$ws.Range("L25").Value = -18.565400843881

# Row 26
$ws.Range("C26").Value = 116
$ws.Range("D26").Value = 92
$ws.Range("E26").Value = 26.086956521739
$ws.Range("F26").Value = 418
$ws.Range("G26").Value = 430
$ws.Range("H26").Value = -2.790697674418
$ws.Range("I26").Value = 275
$ws.Range("J26").Value = 283
$ws.Range("K26").Value = -2.826855123674
$ws.Range("L26").Value = -7.407407407407
$ws.Range("M26").Value = -27.821522309711

# Row 27
$ws.Range("C27").Value = 7
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 23
$ws.Range("G27").Value = 24
$ws.Range("H27").Value = -4.166666666666
$ws.Range("I27").Value = 17
$ws.Range("J27").Value = 19
$ws.Range("K27").Value = -10.526315789473
$ws.Range("L27").Value = -32

# Row 28
$ws.Range("C28").Value = 18
$ws.Range("D28").Value = 6
$ws.Range("E28").Value = 200
$ws.Range("F28").Value = 53
$ws.Range("G28").Value = 35
$ws.Range("H28").Value = 51.428571428571
$ws.Range("I28").Value = 40
$ws.Range("J28").Value = 20
$ws.Range("K28").Value = 100
$ws.Range("L28").Value = 14.285714285714

# Row 29
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = -40
$ws.Range("F29").Value = 12
$ws.Range("G29").Value = 16
$ws.Range("H29").Value = -25
$ws.Range("I29").Value = 10
$ws.Range("J29").Value = 10
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 25
$ws.Range("M29").Value = -47.368421052631
$ws.Range("N29").Value = -90.74074074074

# Row 30
$ws.Range("D30").Value = 4
$ws.Range("E30").Value = -50
$ws.Range("F30").Value = 8
$ws.Range("G30").Value = 13
$ws.Range("H30").Value = -38.461538461538
$ws.Range("I30").Value = 7
$ws.Range("J30").Value = 9
$ws.Range("K30").Value = -22.222222222222
$ws.Range("L30").Value = -12.5
$ws.Range("M30").Value = -53.333333333333
$ws.Range("N30").Value = -92.929292929292

# Row 31
$ws.Range("G31").Value = 9
$ws.Range("J31").Value = 7

# Row 40
$ws.Range("J40").Value = 248
$ws.Range("K40").Value = -32.054794520547
$ws.Range("L40").Value = -45.969498910675
$ws.Range("M40").Value = -58.249158249158
$ws.Range("N40").Value = -64.971751412429

# Row 41
$ws.Range("J41").Value = 2362
$ws.Range("K41").Value = -55.069431234544
$ws.Range("L41").Value = -66.543909348441
$ws.Range("M41").Value = -85.97803502523
$ws.Range("N41").Value = -88.483104978302

# Row 42
$ws.Range("J42").Value = 4257
$ws.Range("K42").Value = -6.82862770847
$ws.Range("L42").Value = -18.463895805401
$ws.Range("M42").Value = -49.261025029797
$ws.Range("N42").Value = -57.336139506915

# Row 43
$ws.Range("J43").Value = 1893
$ws.Range("K43").Value = -55.605065666041
$ws.Range("L43").Value = -65.487693710118
$ws.Range("M43").Value = -84.369581372306
$ws.Range("N43").Value = -87.374107917027

# Row 45
$ws.Range("J45").Value = 1717
$ws.Range("K45").Value = -49.692352768825
$ws.Range("L45").Value = -54.827676927124
$ws.Range("M45").Value = -82.0735017749
$ws.Range("N45").Value = -86.320905035054

# Row 46
$ws.Range("J46").Value = 15847
$ws.Range("K46").Value = -27.87967050471
$ws.Range("L46").Value = -38.696324951644
$ws.Range("M46").Value = -71.145302257829
$ws.Range("N46").Value = -76.70689224347

